$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.075949192047119
$ws.Range("B1").Value = 4.009919166564941
$ws.Range("C1").Value = 2.516944169998169
$ws.Range("D1").Value = 1.868634700775146
$ws.Range("E1").Value = 1.119463801383972
